$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width (in characters) that column M ("In Advance") currently has,
# so the newly inserted column can reuse the same displayed width.
$mWidth = $ws.Columns.Item(13).ColumnWidth

# Insert a new, blank column before column N ("Late"), pushing
# N -> O, O -> P, P -> Q (Late / Paid Date / Outstanding).
$ws.Columns.Item(14).Insert()

# Give the newly inserted (now empty) column N the same width as column M.
$ws.Columns.Item(14).ColumnWidth = $mWidth

# Make "Repayment schedule" the active/selected sheet and reset its view
# (no more frozen/scrolled topLeftCell), selecting cell R6 as the last
# active cell, matching how the sheet was left after the edit.
$ws.Activate()
$ws.Range("R6").Select()
